$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = @'
Background
id="Par1">COVID-19 outbreak has led to severe health burden in the elderly.

 Age, morbidity and dementia have been associated with adverse outcome.


Aims
id="Par2">To evaluate the impact of COVID-19 on health status in home-dwelling patients.


Methods
id="Par3">848 home-dwelling outpatients with dementia contacted from April 27 to 30 and evaluated by a semi-structured interview to evaluate possible health complication due to COVID-19 from February 21 to April 30. Age, sex, education, clinical characteristics (including diagnosis of dementia) and flu vaccination history were obtained from previous medical records.

 Items regarding change in health status and outcome since the onset of the outbreak were collected.

 COVID-19 was diagnosed in patients who developed symptoms according to WHO criteria or tested positive at nasal/throat swab if hospitalized.

 Unplanned hospitalization, institutionalization and mortality were recorded.


Results
id="Par4">Patients were 79.7 years old (SD 7.1) and 63.1% were females.

 Ninety-five (11.2%) patients developed COVID-19-like symptoms.

 Non COVID-19 and COVID-19 patients differed for frequency of diabetes (18.5% vs.

 37.9%, p &lt; 0.001), COPD (7.3% vs.

 18.9%, p &lt; 0.001), and previous flu vaccination (56.7% vs.

 37.9%, p &lt; 0.001).

 Diabetes and COPD were positively associated with COVID-19, whereas higher dementia severity and flu vaccination showed an inverse association.

 Among COVID-19 patients, 42 (44.2%) were hospitalized while 32 (33.7%) died.

 Non COVID-19 patients’ hospitalization and mortality rate were 1.9% and 1.2%, respectively.

 COVID-19 and COPD were significantly associated with the rate of mortality.


Discussion/conclusions
id="Par5">A high proportion of adverse outcome related to COVID-19 was observed in home-dwelling elderly patients with dementia.

 Active monitoring though telehealth programs would be useful particularly for those at highest risk of developing COVID-19 and its adverse outcomes.


Electronic supplementary material
The online version of this article (10.1007/s40520-020-01676-z) contains supplementary material, which is available to authorized users.



'@

$ws.Range("E9").Value = @'
[Salvatore%Caratozzolo%NULL%1,  Alberto%Zucchelli%NULL%1,  Marinella%Turla%NULL%1,  Maria Sofia%Cotelli%NULL%1,  Sara%Fascendini%NULL%1,  Mara%Zanni%NULL%1,  Angelo%Bianchetti%NULL%1,  Matteo Peli%Psy%NULL%1,  Renzo%Rozzini%NULL%1,  Stefano%Boffelli%NULL%1,  Melania%Cappuccio%NULL%1,  Federica Gottardi%Psy%NULL%1,  Chiara Vecchi%Psy%NULL%1,  Daniele%Bellandi%NULL%1,  Claudia%Caminati%NULL%1,  Simona%Gentile%NULL%1,  Elena Lucchi%Psy%NULL%1,  Ignazio%Di Fazio%NULL%1,  Marina Zanetti%Psy%NULL%1,  Giuliana%Vezzadini%NULL%1,  Chiara Forlani%Psy%NULL%1,  Maura Cosseddu%Psy%NULL%1,  Rosanna Turrone%Psy%NULL%1,  Silvia Pelizzari%Psy%NULL%1,  Andrea%Scalvini%NULL%1,  Marco%Di Cesare%NULL%1,  Marta Grigolo%Psy%NULL%1,  Lina%Falanga%NULL%1,  Nives%Medici%NULL%1,  Nives%Palamini%NULL%1,  Elisa Zanacchi%Psy%NULL%1,  Eleonora Grossi%Psy%NULL%1,  Giuseppe%Bellelli%NULL%1,  Alessandra%Marengoni%NULL%1,  Marco%Trabucchi%NULL%1,  Alessandro%Padovani%alessandro.padovani@unibs.it%1,  NULL%NULL%NULL%0]
'@

$ws.Range("D11").Value = @'
Objective
id="Par1">The aim of this study was to investigate the clinical characteristics of Corona Virus Disease 2019 in Taizhou, China.


Methods
id="Par2">A single center retrospective observational study was performed between Jan 1, 2020 and Mar 11, 2020 at Taizhou Public Health Medical Center, Zhejiang, China.

 All patients with confirmed Corona Virus Disease 2019 were enrolled, and their clinical data were gathered by reviewing electronic medical records.

 Outcomes of severely ill patients and non-severely ill patients were compared.


Results
id="Par3">Of 145 hospitalized patients with COVID-19, the average age was 47.5 years old (standard deviation, 14.6) and 54.5% were men.

 Hypertension was the most common comorbidity (15.2%), followed by diabetes mellitus (9.7%).

 Common symptoms included dry cough (81.4%), fever (75.2%), anorexia (42.8%), fatigue (40.7%), chest tightness (32.4%), diarrhea (26.9%) and dizziness (20%).

 According to imaging examination, 79.3% patients showed bilateral pneumonia, 18.6% showed unilateral pneumonia, 61.4% showed ground-glass opacity, and 2.1% showed no abnormal result.

 Compared with non-severely ill patients, severely ill patients were older (mean, years, 52.8 vs.

 45.3, p &lt; 0.01), had a higher proportion of diabetes mellitus (16.3% vs.

 6.9%, p = 0.08), had a higher body mass index (mean, 24.78 vs.

 23.20, p = 0.02) and were more likely to have fever (90.7% vs.

 68.6%, p = 0.01), anorexia (60.5% vs.

 35.3%, p = 0.01), chest tightness (60.5% vs.

20.6%, p &lt; 0.01) and dyspnea (7.0% vs.

 0%, p = 0.03).

 Of the 43 severely ill patients, 6 (14%) received high-flow nasal cannula oxygen therapy, and 1 (2.3%) received invasive mechanical ventilation.


Conclusions
id="Par4">Older patients or patients with comorbidities such as obesity or diabetes mellitus were more likely to have severe condition.

 Treatments of COVID-19 is still experimental and more clinical trials are needed.



'@

$ws.Range("E11").Value = @'
[Qingqing%Chen%NULL%2,  Zhencang%Zheng%NULL%1,  Chao%Zhang%NULL%0,  Xijiang%Zhang%NULL%1,  Huijuan%Wu%NULL%1,  Jingdong%Wang%NULL%1,  Shuwei%Wang%NULL%1,  Cheng%Zheng%dr.zhengcheng@foxmail.com%1]
'@

$ws.Range("D49").Value = @'
Objectives
id="Par1">To characterize the chest computed tomography (CT) findings of severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) according to clinical severity.

 We compared the CT features of common cases and severe cases, symptomatic patients and asymptomatic patients, and febrile and afebrile patients.


Methods
id="Par2">This was a retrospective analysis of the clinical and thoracic CT features of 120 consecutive patients with confirmed SARS-CoV-2 pneumonia admitted to a tertiary university hospital between January 10 and February 10, 2020, in Wuhan city, China.


Results
id="Par3">On admission, the patients generally complained of fever, cough, shortness of breath, and myalgia or fatigue, with diarrhea often present in severe cases.

 Severe patients were 20 years older on average and had comorbidities and an elevated lactate dehydrogenase (LDH) level.

 There were no differences in the CT findings between asymptomatic and symptomatic common type patients or between afebrile and febrile patients, defined according to Chinese National Health Commission guidelines.


Conclusions
id="Par4">The clinical and CT features at admission may enable clinicians to promptly evaluate the prognosis of patients with SARS-CoV-2 pneumonia.

 Clinicians should be aware that clinically silent cases may present with CT features similar to those of symptomatic common patients.


Key Points
id="Par5">
• The clinical features and predominant patterns of abnormalities on CT for asymptomatic, typic common, and severe cases were summarized.

 These findings may help clinicians to identify severe patients quickly at admission.



Electronic supplementary material
The online version of this article (10.1007/s00330-020-06854-1) contains supplementary material, which is available to authorized users.



'@

$ws.Range("E49").Value = @'
[Rui%Zhang%NULL%2,  Huangqing%Ouyang%NULL%2,  Lingli%Fu%NULL%2,  Shijie%Wang%NULL%2,  Jianglong%Han%NULL%2,  Kejie%Huang%NULL%2,  Mingfang%Jia%NULL%2,  Qibin%Song%NULL%2,  Zhenming%Fu%davidfuzming@whu.edu.cn%2]
'@
